$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "Leonardo"
$ws.Range("B10").Value = "Fragoso"
$ws.Range("C10").Value = "LeoFragoso"
$ws.Range("D10").Value = "'21980292791"
$ws.Range("E10").Value = "leonardorfragoso@gmail.com"
$ws.Range("F10").Value = "'2"
$ws.Range("G10").Value = "Igual este"
